$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.745.19"
$ws.Range("E2").Value = "  +0.30%  "

$ws.Range("D3").Value = "2.408.33"
$ws.Range("E3").Value = "  -3.14%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'488.68"
$ws.Range("E5").Value = "  -0.82%  "

$ws.Range("D6").Value = "'153.18"
$ws.Range("E6").Value = "  +0.14%  "

$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("D8").Value = "'0.604"
$ws.Range("E8").Value = "  +18.17%  "

$ws.Range("D9").Value = "2.419.14"
$ws.Range("E9").Value = "  -3.04%  "

$ws.Range("D10").Value = "'6.11"
$ws.Range("E10").Value = "  +6.19%  "

$ws.Range("D11").Value = "'0.0999"
$ws.Range("E11").Value = "  +1.48%  "

$ws.Range("D12").Value = "'0.335"
$ws.Range("E12").Value = "  +0.39%  "

$ws.Range("D13").Value = "'0.126"
$ws.Range("E13").Value = "  +1.32%  "

$ws.Range("D14").Value = "2.825.88"
$ws.Range("E14").Value = "  -3.05%  "

$ws.Range("D15").Value = "56.886.88"
$ws.Range("E15").Value = "  +0.23%  "

$ws.Range("D16").Value = "'20.81"
$ws.Range("E16").Value = "  -2.51%  "

$ws.Range("E17").Value = "  -2.09%  "

$ws.Range("D18").Value = "2.420.15"
$ws.Range("E18").Value = "  -3.01%  "

$ws.Range("D19").Value = "'4.76"
$ws.Range("E19").Value = "  +4.51%  "

$ws.Range("D20").Value = "'324.77"
$ws.Range("E20").Value = "  +1.35%  "

$ws.Range("D21").Value = "'10.01"
$ws.Range("E21").Value = "  -2.87%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'5.99"
$ws.Range("E22").Value = "  +2.00%  "

$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.15%  "

$ws.Range("D24").Value = "'58.07"
$ws.Range("E24").Value = "  -1.28%  "

$ws.Range("D25").Value = "'0.409"
$ws.Range("E25").Value = "  -0.21%  "

$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.15%  "

$ws.Range("E27").Value = "  -1.26%  "

$ws.Range("D28").Value = "2.518.58"
$ws.Range("E28").Value = "  -3.36%  "

$ws.Range("D29").Value = "'7.33"
$ws.Range("E29").Value = "  -3.44%  "

$ws.Range("D30").Value = "0.0₃0784"
$ws.Range("E30").Value = "  -3.04%  "

$ws.Range("E31").Value = "  +0.12%  "

$ws.Range("D32").Value = "'150.07"
$ws.Range("E32").Value = "  -0.83%  "

$ws.Range("D33").Value = "'18.57"
$ws.Range("E33").Value = "  +1.33%  "

$ws.Range("E34").Value = "  +0.38%  "

$ws.Range("D35").Value = "'5.33"
$ws.Range("E35").Value = "  +1.24%  "

$ws.Range("D36").Value = "'1.16"
$ws.Range("E36").Value = "  -0.36%  "

$ws.Range("D37").Value = "'3.73"
$ws.Range("E37").Value = "  -1.35%  "

$ws.Range("D38").Value = "'0.847"
$ws.Range("E38").Value = "  -1.79%  "

$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").Value = "'3.54"
$ws.Range("E39").Value = "  +0.69%  "

$ws.Range("D40").Value = "'34.06"
$ws.Range("E40").Value = "  +0.31%  "

$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").Value = "'0.101"
$ws.Range("E41").Value = "  +9.25%  "

$ws.Range("D42").Value = "'1.37"
$ws.Range("E42").Value = "  -1.54%  "

$ws.Range("D43").Value = "'0.994"
$ws.Range("E43").Value = "  -0.20%  "

$ws.Range("D44").Value = "'0.593"
$ws.Range("E44").Value = "  -3.69%  "

$ws.Range("D45").Value = "'268.15"
$ws.Range("E45").Value = "  +1.26%  "

$ws.Range("D46").Value = "'0.0534"
$ws.Range("E46").Value = "  -5.11%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0229"
$ws.Range("E47").Value = "  -0.26%  "

$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D48").Value = "'10.20"
$ws.Range("E48").Value = "  -0.18%  "

$ws.Range("D49").Value = "'4.63"
$ws.Range("E49").Value = "  -6.74%  "

$ws.Range("D50").Value = "'17.44"
$ws.Range("E50").Value = "  -1.90%  "

$ws.Range("D51").Value = "1.867.47"
$ws.Range("E51").Value = "  -1.10%  "
